$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 543
$ws1.Range("F7").Value  = 101
$ws1.Range("F8").Value  = 116
$ws1.Range("F9").Value  = 43
$ws1.Range("F10").Value = 6690
$ws1.Range("F11").Value = 231
$ws1.Range("F12").Value = 364
$ws1.Range("F13").Value = 2954
$ws1.Range("F14").Value = 191
$ws1.Range("F15").Value = 329
$ws1.Range("F17").Value = 535

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 543
$ws4.Range("F9").Value  = 101
$ws4.Range("F10").Value = 116
$ws4.Range("F11").Value = 43
$ws4.Range("F13").Value = 6690
$ws4.Range("F15").Value = 231
$ws4.Range("F16").Value = 364
$ws4.Range("F17").Value = 2954
$ws4.Range("F18").Value = 191
$ws4.Range("F19").Value = 329
$ws4.Range("F21").Value = 535
